$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "manualStatus" column (I) for rows 2-13 from the numeric
# value 128 to the text value "[128]".
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = "[128]"
}

# Move the active selection to I13, matching the saved selection state.
$ws.Range("I13").Select()
